$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (표준내역 구분): "주택(순타)" -> "주택" for data rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = "주택"
}

# Update column H (소속부서): "건적팀" -> "견적팀" for data rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = "견적팀"
}
